# ReactMindMap.xlsx edit script
# Commit message: "upto 9 list and keys"
# Adds "Handling Events", "Conditional Rendering" and "Lists and Keys" sections
# to the React-MainConcepts sheet, renames that sheet to React-16.5.2-MainConcepts,
# adds a defined name "ReactDocs" pointing at the renamed sheet, and adds two
# external hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("React-MainConcepts")
$ws.Activate()

# ---------------------------------------------------------------------------
# New shared-string text used below (kept in variables for readability)
# ---------------------------------------------------------------------------
$s76 = "Handling Events"
$s77 = "React events are named using camelCase, rather than lowercase."
$s78 = "With JSX you pass a function as the event handler, rather than a string."
$s79 = "<button onClick={activateLasers}> Activate Lasers</button>"
$s80 = " cannot return false to prevent default behavior in React. You must call preventDefault explicitly."
$s81 = "event obj passed to handler is synthetic event"
$s82 = "Conditional Rendering"
$s83 = "Preventing Component from Rendering. Return null from render() method of component"
$s84 = "can render component in if/else or for loop"
$s85 = "Returning null from a component’s render method does not affect the firing of the component’s lifecycle methods. For instance componentDidUpdate will still be called."
$s86 = "Lists and Keys"
$s87 = "Keys"
$s88 = "Keys help React identify which items have changed, are added, or are removed"
$s89 = "Keys give the array elements a stable identity"
$s90 = "not recommend using indexes for keys if the order of items may change."
$s91 = "https://medium.com/@robinpokorny/index-as-a-key-is-an-anti-pattern-e0349aece318"
$s92 = "https://reactjs.org/docs/reconciliation.html#recursing-on-children"
$s93 = "A good rule of thumb is that elements inside the map() call need keys."
$s94 = "identifies a list item among its siblings. "
$s95 = "We can use the same keys when we produce two different arrays"
$s96 = "Keys serve as a hint to React but they don’t get passed to your components. If you need the same value in your component, pass it explicitly as a prop with a different name"
$s97 = "const content = posts.map((post) =>`n  <Post`n    key={post.id}`n    id={post.id}`n    title={post.title} />`n);"
$s98 = "Keep in mind that if the map() body is too nested, it might be a good time to extract a component."

# ---------------------------------------------------------------------------
# Row 38-40: Handling Events
# ---------------------------------------------------------------------------
$ws.Range("B38").Value = $s76
$ws.Range("B38").Font.Bold = $true

$ws.Range("C38").Value = $s77
$ws.Range("C38").WrapText = $true

$ws.Range("D38").Value = $s79
$ws.Range("D38:D39").Merge()
$ws.Range("D38:D39").VerticalAlignment = -4160

$ws.Range("C39").Value = $s78

$ws.Range("C40").Value = $s80
$ws.Range("D40").Value = $s81

# ---------------------------------------------------------------------------
# Row 42-44: Conditional Rendering
# ---------------------------------------------------------------------------
$ws.Range("B42").Value = $s82
$ws.Range("B42").Characters(1, 11).Font.Bold = $true
$ws.Range("B42").Characters(13, 9).Font.Bold = $true

$ws.Range("C42").Value = $s83
$ws.Range("C43").Value = $s84

$ws.Range("C44").Value = $s85
$ws.Range("C44").WrapText = $true
$ws.Rows.Item(44).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 46-54: Lists and Keys
# ---------------------------------------------------------------------------
$ws.Range("B46").Value = $s86
$ws.Range("B46").Font.Bold = $true

$ws.Range("C46").Value = $s87
$ws.Range("D46").Value = $s92
$ws.Hyperlinks.Add($ws.Range("D46"), "https://reactjs.org/docs/reconciliation.html", "recursing-on-children")

$ws.Range("C47").Value = $s88
$ws.Range("C48").Value = $s89
$ws.Range("C49").Value = $s94

$ws.Range("C50").Value = $s95
$ws.Range("D50").Value = $s91
$ws.Hyperlinks.Add($ws.Range("D50"), $s91)

$ws.Range("C51").Value = $s90
$ws.Range("C52").Value = $s93

$ws.Range("C53").Value = $s96
$ws.Range("C53").WrapText = $true
$ws.Range("C53").VerticalAlignment = -4160
$ws.Rows.Item(53).RowHeight = 90

$ws.Range("D53").Value = $s97
$ws.Range("D53").WrapText = $true

$ws.Range("C54").Value = $s98

# ---------------------------------------------------------------------------
# Rename the sheet and add a defined name referencing it
# ---------------------------------------------------------------------------
$ws.Name = "React-16.5.2-MainConcepts"
$wb.Names.Add("ReactDocs", "='React-16.5.2-MainConcepts'!`$1:`$1048576")

# ---------------------------------------------------------------------------
# Update the view: select the full-row range (mirrors the state captured
# when the defined name's range was last selected in the source file).
# ---------------------------------------------------------------------------
$ws.Rows("1:1048576").Select()
